$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '95.036.22'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.475.18'
$ws.Range('E3').Value = '  +4.41%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.35'
$ws.Range('E5').Value = '  -3.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '645.15'
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E7').Value = '  +6.70%  '
$ws.Range('E8').Value = '  -3.17%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +2.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.474.76'
$ws.Range('E11').Value = '  +4.49%  '
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.30'
$ws.Range('E12').Value = '  +5.15%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.199'
$ws.Range('E13').Value = '  -3.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.16'
$ws.Range('E14').Value = '  +1.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.934.18'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.127.95'
$ws.Range('E16').Value = '  +4.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000257'
$ws.Range('E17').Value = '  +2.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.54'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.477.04'
$ws.Range('E19').Value = '  +4.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.00'
$ws.Range('E20').Value = '  +5.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.47'
$ws.Range('E21').Value = '  +9.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.515'
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '503.03'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  -5.33%  '
$ws.Range('E25').Value = '  -1.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.51'
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '92.06'
$ws.Range('E27').Value = '  -3.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.19'
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.661.39'
$ws.Range('E29').Value = '  +4.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '11.79'
$ws.Range('E30').Value = '  +7.61%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('E32').Value = '  +12.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.185'
$ws.Range('E34').Value = '  -1.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '30.97'
$ws.Range('E35').Value = '  +11.17%  '
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.570'
$ws.Range('E37').Value = '  +4.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.83'
$ws.Range('E38').Value = '  +3.17%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '537.02'
$ws.Range('E39').Value = '  +6.14%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.45'
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.932'
$ws.Range('E41').Value = '  +12.74%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '24.10'
$ws.Range('E44').Value = '  -0.97%  '
$ws.Range('E45').Value = '  +4.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.71'
$ws.Range('E46').Value = '  +2.60%  '
$ws.Range('E47').Value = '  -2.52%  '
$ws.Range('E48').Value = '  -3.21%  '
$ws.Range('E49').Value = '  +10.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.23'
$ws.Range('E50').Value = '  +4.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.36'
$ws.Range('E51').Value = '  -0.07%  '
